{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"List Bullet\") {\n    para.insertText(\"Design: \", Word.InsertLocation.start);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"List Bullet\") {\n        $p.Range.InsertBefore(\"Design: \")\n    }\n}\n"}
